# Replace the categorical "color name" placeholder values in column B
# (rows 2-11) of each color sheet with the actual numeric data values.
# The "Blue" sheet never received real data, so its B2:B11 placeholder
# cells are simply cleared out instead of being populated with numbers.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> array of 10 numeric values for rows 2..11
# (Toyota Yaris, Mazda MX-30, Honda JAZZ, Land Rover Defender, SEAT Leon,
#  KIA Sorento, Honda e, Hyundai i10, ISUZU D-Max Crew Cab, Audi A3)
$data = @{
    "Default green" = @(18, 0, 27, 0, 0, 0, 23, 34, 0, 3)
    "Green"         = @(53, 78, 36, 78, 49, 41, 56, 13, 23, 37)
    "Yellow"        = @(31.5, 45.75, 30, 20.25, 30.75, 49.5, 22.5, 24.75, 63, 31.5)
    "Orange"        = @(7.5, 4.5, 9.5, 17.5, 16, 15, 7.5, 18, 18.5, 18)
    "Brown"         = @(1, 3, 2.25, 2.5, 3.75, 4.5, 1.5, 3.5, 3.5, 5.25)
    "Red"           = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "Default Red"   = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $data[$sheetName]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }
}

# The "Blue" sheet's placeholder text values are removed entirely (no data
# was ever collected for this color), leaving the data cells blank.
$wsBlue = $wb.Worksheets.Item("Blue")
$wsBlue.Range("B2:B11").ClearContents()
